$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook calculation settings: enable iterative calculation ---
$excel.Iteration = $true
$excel.MaxIterations = 1000
$excel.CalculateBeforeSave = $false

# --- Update data values (columns D:O get recalculated/updated figures, column P is the new 2022 year) ---
$ws.Range("P4").Value = 2022
$ws.Range("F5").Value = 61.011419249592166
$ws.Range("G5").Value = 55.628058727569339
$ws.Range("H5").Value = 57.748776508972263
$ws.Range("I5").Value = 66.068515497553022
$ws.Range("J5").Value = 64.763458401305058
$ws.Range("K5").Value = 64.600326264274059
$ws.Range("L5").Value = 76.508972267536706
$ws.Range("M5").Value = 71.125611745513879
$ws.Range("N5").Value = 70.96247960848288
$ws.Range("O5").Value = 72.920065252854812
$ws.Range("P5").Value = 76.508972267536706
$ws.Range("F6").Value = 107.1
$ws.Range("G6").Value = 109.25
$ws.Range("H6").Value = 109.89999999999999
$ws.Range("I6").Value = 115.75
$ws.Range("J6").Value = 115.14999999999999
$ws.Range("K6").Value = 114.8
$ws.Range("L6").Value = 112.15000000000002
$ws.Range("M6").Value = 110.6
$ws.Range("N6").Value = 109.3
$ws.Range("O6").Value = 110.94999999999999
$ws.Range("P6").Value = 110.15000000000002
$ws.Range("F7").Value = 44.821917808219176
$ws.Range("G7").Value = 46.246575342465754
$ws.Range("H7").Value = 45.643835616438352
$ws.Range("I7").Value = 48.328767123287676
$ws.Range("J7").Value = 44.986301369863014
$ws.Range("K7").Value = 48.493150684931507
$ws.Range("L7").Value = 48.876712328767127
$ws.Range("M7").Value = 53.424657534246577
$ws.Range("N7").Value = 55.178082191780817
$ws.Range("O7").Value = 56.821917808219183
$ws.Range("P7").Value = 66.630136986301366
$ws.Range("F8").Value = 79.296875
$ws.Range("G8").Value = 83.203125
$ws.Range("H8").Value = 78.515625
$ws.Range("I8").Value = 73.4375
$ws.Range("J8").Value = 88.28125
$ws.Range("K8").Value = 91.406249999999986
$ws.Range("L8").Value = 99.609375
$ws.Range("M8").Value = 98.437499999999986
$ws.Range("N8").Value = 62.890625
$ws.Range("O8").Value = 92.578124999999986
$ws.Range("P8").Value = 102.34375
$ws.Range("F9").Value = 119.38663745892661
$ws.Range("G9").Value = 139.10186199342826
$ws.Range("H9").Value = 141.29244249726176
$ws.Range("I9").Value = 148.95947426067906
$ws.Range("J9").Value = 119.38663745892661
$ws.Range("K9").Value = 116.10076670317633
$ws.Range("L9").Value = 124.8630887185104
$ws.Range("M9").Value = 134.72070098576123
$ws.Range("N9").Value = 166.48411829134719
$ws.Range("O9").Value = 167.57940854326395
$ws.Range("P9").Value = 187.29463307776561
$ws.Range("F10").Value = 101.11561866125763
$ws.Range("G10").Value = 101.52129817444219
$ws.Range("H10").Value = 101.82555780933065
$ws.Range("I10").Value = 111.15618661257606
$ws.Range("J10").Value = 111.56186612576064
$ws.Range("K10").Value = 111.25760649087222
$ws.Range("L10").Value = 146.45030425963489
$ws.Range("M10").Value = 128.39756592292088
$ws.Range("N10").Value = 107.80933062880325
$ws.Range("P10").Value = 107.20081135902637
$ws.Range("F11").Value = 132.37095363079615
$ws.Range("G11").Value = 132.45844269466318
$ws.Range("H11").Value = 133.59580052493439
$ws.Range("I11").Value = 142.25721784776903
$ws.Range("J11").Value = 144.61942257217851
$ws.Range("K11").Value = 156.95538057742783
$ws.Range("L11").Value = 163.95450568678916
$ws.Range("M11").Value = 143.48206474190727
$ws.Range("N11").Value = 161.41732283464566
$ws.Range("O11").Value = 152.1434820647419
$ws.Range("P11").Value = 155.38057742782152
$ws.Range("F12").Value = 19.563459983831848
$ws.Range("G12").Value = 23.686337914308812
$ws.Range("H12").Value = 24.00970088924818
$ws.Range("I12").Value = 26.434923201293454
$ws.Range("J12").Value = 35.89329021827001
$ws.Range("K12").Value = 35.569927243330639
$ws.Range("L12").Value = 27.081649151172186
$ws.Range("M12").Value = 33.87227162489895
$ws.Range("N12").Value = 26.919967663702504
$ws.Range("O12").Value = 26.273241713823765
$ws.Range("P12").Value = 36.054971705739689
$ws.Range("F13").Value = 147.56902027852431
$ws.Range("G13").Value = 150.25653554849742
$ws.Range("H13").Value = 150.37869533349621
$ws.Range("I13").Value = 149.03493769850965
$ws.Range("J13").Value = 147.44686049352552
$ws.Range("K13").Value = 133.64280478866357
$ws.Range("L13").Value = 102.12558025897873
$ws.Range("M13").Value = 99.315905204006839
$ws.Range("N13").Value = 107.37845101392622
$ws.Range("O13").Value = 85.023210359149758
$ws.Range("P13").Value = 106.76765208893232
$ws.Range("F14").Value = 21.978021978021978
$ws.Range("J14").Value = 10.989010989010989
$ws.Range("K14").Value = 10.989010989010989
$ws.Range("L14").Value = 12.087912087912089
$ws.Range("M14").Value = 6.593406593406594
$ws.Range("N14").Value = 10.989010989010989
$ws.Range("P14").Value = 32.967032967032964


# --- Copy cell formatting (style) from the adjacent 2021 column (O) into the new 2022 column (P) ---
$ws.Range("O3").Copy()
$ws.Range("P3").PasteSpecial(-4122)
$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)
$ws.Range("O5").Copy()
$ws.Range("P5").PasteSpecial(-4122)
$ws.Range("O6").Copy()
$ws.Range("P6").PasteSpecial(-4122)
$ws.Range("O7").Copy()
$ws.Range("P7").PasteSpecial(-4122)
$ws.Range("O8").Copy()
$ws.Range("P8").PasteSpecial(-4122)
$ws.Range("O9").Copy()
$ws.Range("P9").PasteSpecial(-4122)
$ws.Range("O10").Copy()
$ws.Range("P10").PasteSpecial(-4122)
$ws.Range("O11").Copy()
$ws.Range("P11").PasteSpecial(-4122)
$ws.Range("O12").Copy()
$ws.Range("P12").PasteSpecial(-4122)
$ws.Range("O13").Copy()
$ws.Range("P13").PasteSpecial(-4122)
$ws.Range("O14").Copy()
$ws.Range("P14").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Move the active selection to R1, matching the saved view state ---
$ws.Range("R1").Select()
